$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell H1 ("Save") should carry the same style as the other
# header cells (bold / bordered / centered -> style index 1 in the source).
# Copy the format from the existing "sum" header (G1) onto H1, then set text.
$ws.Range("G1").Copy() | Out-Null
$ws.Range("H1").PasteSpecial(-4122) | Out-Null
$ws.Range("H1").Value = "Save"

# New "Save" data column values (H2:H10), plain numbers, no special style.
$saveValues = @(1, 0, 0, 0, 0, 0, 0, 0, 0)
for ($i = 0; $i -lt $saveValues.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 8).Value = $saveValues[$i]
}
